$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("N2mD_test")
$ws.Activate()

# ---------------------------------------------------------------------------
# The reaction equation for the EC_biomass row (C18) simplified from a
# weighted multi-substrate reaction to a single-substrate reaction.
# (Edited before the header text below so new shared-string entries land
# in the same order as the authored workbook.)
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = "[c] : B ---> F"

# ---------------------------------------------------------------------------
# Header row: units on the Metab low/high columns changed from mM to M
# ---------------------------------------------------------------------------
$ws.Range("T1").Value = "Mclow (M)"
$ws.Range("U1").Value = "Mchigh (M)"

# ---------------------------------------------------------------------------
# Regulator Km values (columns T/U, rows 2-8): tightened from mM-scale
# (0.1 / 20) placeholder values down to M-scale (1E-4 / 0.02) values.
# ---------------------------------------------------------------------------
$rows = 2..8
foreach ($r in $rows) {
    $ws.Cells.Item($r, 20).Value = 0.0001   # column T
    $ws.Cells.Item($r, 21).Value = 0.02     # column U
}

# ---------------------------------------------------------------------------
# Metab low/high (column M) updates for the exchange / biomass rows that
# previously held placeholder zeros.
# ---------------------------------------------------------------------------
$ws.Range("M12").Value = 0.33
$ws.Range("M14").Value = 5
$ws.Range("M15").Value = 2.5
$ws.Range("M17").Value = 0.33
$ws.Range("M18").Value = 0.33

# ---------------------------------------------------------------------------
# Restore the active cell selection to what it was when the author saved.
# ---------------------------------------------------------------------------
$ws.Range("U4").Select()
